# The data rows (2-32) of the sheet have been reshuffled: each destination
# row now carries the values (Fecha, Variedad, Calidad, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades) that used to belong to a different row, while columns
# A, B, C, E, F, G, R stay the same for every row. The mapping below says,
# for each destination row, which original row's values should be copied in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2 = 4;  3 = 13; 4 = 32; 5 = 17; 6 = 8;  7 = 26; 8 = 9;  9 = 2;  10 = 29;
    11 = 10; 12 = 30; 13 = 20; 14 = 21; 15 = 22; 16 = 28; 17 = 25; 18 = 31;
    19 = 6; 20 = 7; 21 = 14; 22 = 23; 23 = 5; 24 = 27; 25 = 3; 26 = 11;
    27 = 24; 28 = 16; 29 = 15; 30 = 18; 31 = 19; 32 = 12
}

$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the original values of every row before any writes happen, so
# that later writes don't clobber data that still needs to be read.
$snapshot = @{}
for ($r = 2; $r -le 32; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
